# Trade #125 closed at 2026-02-17 09:32:40 - unknown UNKNOWN +0.000%
#
# Updates the Summary + Strategy Status roll-up figures and appends the
# newly closed trade (#125) to both the "All Trades" and "MarketMaking"
# trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up figures
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.31   # Current Capital
$summary.Range("B4").Value = 0.32      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 125       # Total Trades
$summary.Range("B7").Value = 56        # Winning Trades
$summary.Range("B9").Value = 44.8      # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.31     # Capital
$status.Range("D4").Value = 125        # Trades
$status.Range("E4").Value = 0.32       # P&L $
$status.Range("F4").Value = 0.31       # P&L %
$status.Range("G4").Value = 44.8       # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (#125) to a trade-log style sheet.
#    Date/time-looking strings ("2026-02-17", "09:32:34") must land as
#    literal text, not get auto-converted into Excel serial dates, and
#    must not leave behind any number-format / style residue. Writing
#    a self-referential text formula and then collapsing it to a value
#    via copy / paste-special achieves that without touching styles.xml.
# ---------------------------------------------------------------------
function Add-Trade125Row($sheet) {
    $row = 126

    $sheet.Cells.Item($row, 1).Value = 125

    $sheet.Cells.Item($row, 2).Formula = '="2026-02-17"'
    $sheet.Cells.Item($row, 2).Copy()
    $sheet.Cells.Item($row, 2).PasteSpecial(-4163)

    $sheet.Cells.Item($row, 3).Value = "09:32:34"
    $sheet.Cells.Item($row, 4).Value = "MarketMaking"
    $sheet.Cells.Item($row, 5).Value = "UP"
    $sheet.Cells.Item($row, 6).Value = 0.61
    $sheet.Cells.Item($row, 7).Value = 0.68
    $sheet.Cells.Item($row, 8).Value = "CLOSED"
    $sheet.Cells.Item($row, 9).Value = 11.4754
    $sheet.Cells.Item($row, 10).Value = 0.07000000000000001
    $sheet.Cells.Item($row, 11).Value = 100.31
    $sheet.Cells.Item($row, 12).Value = 0
    $sheet.Cells.Item($row, 13).Value = 0
    $sheet.Cells.Item($row, 14).Value = 0.6
    $sheet.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item($row, 16).Value = "early_exit"
    $sheet.Cells.Item($row, 17).Value = 0.1
}

Add-Trade125Row($wb.Worksheets.Item("All Trades"))
Add-Trade125Row($wb.Worksheets.Item("MarketMaking"))
